# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing worker rows (17, 18) - they shift rows 23/24 up to 21/22.
$ws.Rows("17:18").Delete()

# The remaining single worker row (16) now represents the updated record.
$ws.Range("C16").Value = "1018492269"
$ws.Range("D16").Value = "VICTOR RAFAEL SANJUAN CASSIANI"
$ws.Range("E16").Value = "2508"

# Update the summary total "VALOR MORA" (now only one worker in arrears).
$ws.Range("E11").Value = 56940

# Update the worker count summary.
$ws.Range("C13").Value = 1
